$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cells value while preserving it as literal text
# (so numeric-looking strings like "1.00" or "0.0610" keep their exact
# formatting, matching the original inlineStr cells) without leaving the
# cells style/number-format changed afterwards.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "37.058.79"
Set-TextValue $ws.Range("E2") "  -0.66%  "
Set-TextValue $ws.Range("D3") "2.009.04"
Set-TextValue $ws.Range("E3") "  -2.06%  "
Set-TextValue $ws.Range("D4") "1.01"
Set-TextValue $ws.Range("E4") "  +0.96%  "
Set-TextValue $ws.Range("D5") "226.15"
Set-TextValue $ws.Range("E5") "  -1.93%  "
Set-TextValue $ws.Range("E6") "  -2.23%  "
Set-TextValue $ws.Range("E7") "  +0.03%  "
Set-TextValue $ws.Range("D8") "54.54"
Set-TextValue $ws.Range("E8") "  -4.31%  "
Set-TextValue $ws.Range("D9") "0.378"
Set-TextValue $ws.Range("E9") "  -1.34%  "
Set-TextValue $ws.Range("E10") "  +2.31%  "
Set-TextValue $ws.Range("E11") "  -2.96%  "
Set-TextValue $ws.Range("D12") "2.307.27"
Set-TextValue $ws.Range("E12") "  -1.88%  "
Set-TextValue $ws.Range("D13") "14.18"
Set-TextValue $ws.Range("E13") "  -3.22%  "
Set-TextValue $ws.Range("D14") "20.20"
Set-TextValue $ws.Range("E14") "  -1.81%  "
Set-TextValue $ws.Range("E15") "  -2.24%  "
Set-TextValue $ws.Range("E16") "  -2.74%  "
Set-TextValue $ws.Range("D17") "2.023.92"
Set-TextValue $ws.Range("E17") "  -1.47%  "
Set-TextValue $ws.Range("D18") "36.928.66"
Set-TextValue $ws.Range("E18") "  -0.86%  "
Set-TextValue $ws.Range("D19") "6.09"
Set-TextValue $ws.Range("E19") "  +0.86%  "
Set-TextValue $ws.Range("D20") "68.66"
Set-TextValue $ws.Range("E20") "  -1.44%  "
Set-TextValue $ws.Range("D21") "0.0₃0816"
Set-TextValue $ws.Range("E21") "  -0.58%  "
Set-TextValue $ws.Range("D22") "223.18"
Set-TextValue $ws.Range("E22") "  -1.40%  "
Set-TextValue $ws.Range("D23") "1.00"
Set-TextValue $ws.Range("E23") "  -0.01%  "
Set-TextValue $ws.Range("D24") "2.42"
Set-TextValue $ws.Range("E24") "  +1.39%  "
Set-TextValue $ws.Range("E25") "  -5.81%  "
Set-TextValue $ws.Range("D26") "165.11"
Set-TextValue $ws.Range("E26") "  -2.85%  "
Set-TextValue $ws.Range("D27") "9.16"
Set-TextValue $ws.Range("E27") "  -6.22%  "
Set-TextValue $ws.Range("E28") "  -3.61%  "
Set-TextValue $ws.Range("E29") "  +0.38%  "
Set-TextValue $ws.Range("D30") "18.64"
Set-TextValue $ws.Range("E30") "  -2.71%  "
Set-TextValue $ws.Range("E31") "  -3.61%  "
Set-TextValue $ws.Range("D32") "4.51"
Set-TextValue $ws.Range("E32") "  -0.28%  "
Set-TextValue $ws.Range("D33") "0.0610"
Set-TextValue $ws.Range("E33") "  -1.90%  "
Set-TextValue $ws.Range("E34") "  -3.08%  "
Set-TextValue $ws.Range("E35") "  -5.95%  "
Set-TextValue $ws.Range("E36") "  +1.67%  "
Set-TextValue $ws.Range("E37") "  +0.22%  "
Set-TextValue $ws.Range("E38") "  -4.07%  "
Set-TextValue $ws.Range("D39") "5.36"
Set-TextValue $ws.Range("E39") "  +0.91%  "
Set-TextValue $ws.Range("D40") "1.474.18"
Set-TextValue $ws.Range("E40") "  -0.47%  "
Set-TextValue $ws.Range("D41") "0.0215"
Set-TextValue $ws.Range("E41") "  -4.34%  "
Set-TextValue $ws.Range("B42") "InjectiveProtocol"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D42") "16.53"
Set-TextValue $ws.Range("E42") "  +0.37%  "
Set-TextValue $ws.Range("B43") "Aave"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D43") "94.70"
Set-TextValue $ws.Range("E43") "  -3.61%  "
Set-TextValue $ws.Range("D44") "0.0920"
Set-TextValue $ws.Range("E44") "  -3.21%  "
Set-TextValue $ws.Range("E45") "  -5.12%  "
Set-TextValue $ws.Range("D46") "1.12"
Set-TextValue $ws.Range("E46") "  -4.43%  "
Set-TextValue $ws.Range("D47") "7.20"
Set-TextValue $ws.Range("E47") "  -0.52%  "
Set-TextValue $ws.Range("E48") "  -1.99%  "
Set-TextValue $ws.Range("D49") "2.91"
Set-TextValue $ws.Range("E49") "  -1.10%  "
Set-TextValue $ws.Range("D50") "2.199.34"
Set-TextValue $ws.Range("E50") "  -1.69%  "
Set-TextValue $ws.Range("D51") "44.28"
Set-TextValue $ws.Range("E51") "  -2.34%  "
